$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 8 (Connector M/F 01X03 P2.50mm, designator J1) ---
# Designator now also covers the newly added J5 connector, so quantity doubles.
$ws.Range("B8").Value = "J1, J5"
$ws.Range("C8").Value = 4

# --- Insert two new rows before the old last data row ("Power button") ---
$ws.Rows("19:20").Insert()

# Copy formatting (fill/border/number format) from the row above down into the
# two freshly inserted blank rows so they look consistent with the rest of the table.
$ws.Range("A18:F18").Copy()
$ws.Range("A19:F20").PasteSpecial(-4122)

# New row 19: Schottky Diode 1A (D1)
$ws.Range("A19").Value = "Schottky Diode 1A"
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = 0.15
$ws.Range("E19").Formula = "=D19*C19"
$ws.Range("F19").Value = "https://tpetrov.com/bat54j-sod323-40216"
$ws.Range("B19").Value = "D1"

# New row 20: Polyfuse 300 mA (F1)
$ws.Range("A20").Value = "Polyfuse 300 mA"
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 0.61
$ws.Range("E20").Formula = "=C20*D20"
$ws.Range("F20").Value = "https://tpetrov.com/predpazitel-rxef030-244588"
$ws.Range("B20").Value = "F1"

# Update the current selection to match the edited workbook state.
[void]$ws.Range("B24").Select()
